$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the typo in the name cell (row 6, column A)
$ws.Range("A6").Value = "Krishna Sapkota"

# Move the active selection to A6
$ws.Range("A6").Select()
